$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 27, shifting existing rows 27-39
# down to become rows 29-41 (weekly update: two new weeks of data added on top).
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Insert()

# --- New row 27 ---
$ws.Cells.Item(27,1).Value = 8
$ws.Cells.Item(27,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(27,3).Value = "Coquimbo"
$ws.Cells.Item(27,4).Value = 44673
$ws.Cells.Item(27,5).Value = 4
$ws.Cells.Item(27,6).Value = "Fruta"
$ws.Cells.Item(27,7).Value = 100104
$ws.Cells.Item(27,8).Value = "Frutos de pepita"
$ws.Cells.Item(27,9).Value = 100104003
$ws.Cells.Item(27,10).Value = "Membrillo"
$ws.Cells.Item(27,11).Value = "Champion"
$ws.Cells.Item(27,12).Value = "Especial"
$ws.Cells.Item(27,13).Value = 16
$ws.Cells.Item(27,14).Value = 320000
$ws.Cells.Item(27,15).Value = 330000
$ws.Cells.Item(27,16).Value = 325000
$ws.Cells.Item(27,17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(27,18).Value = "Región de O'Higgins"
$ws.Cells.Item(27,19).Value = 722
$ws.Cells.Item(27,20).Value = 450

# --- New row 28 ---
$ws.Cells.Item(28,1).Value = 8
$ws.Cells.Item(28,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28,3).Value = "Coquimbo"
$ws.Cells.Item(28,4).Value = 44673
$ws.Cells.Item(28,5).Value = 4
$ws.Cells.Item(28,6).Value = "Fruta"
$ws.Cells.Item(28,7).Value = 100104
$ws.Cells.Item(28,8).Value = "Frutos de pepita"
$ws.Cells.Item(28,9).Value = 100104003
$ws.Cells.Item(28,10).Value = "Membrillo"
$ws.Cells.Item(28,11).Value = "Champion"
$ws.Cells.Item(28,12).Value = "Primera"
$ws.Cells.Item(28,13).Value = 16
$ws.Cells.Item(28,14).Value = 280000
$ws.Cells.Item(28,15).Value = 290000
$ws.Cells.Item(28,16).Value = 285000
$ws.Cells.Item(28,17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(28,18).Value = "Región de O'Higgins"
$ws.Cells.Item(28,19).Value = 633
$ws.Cells.Item(28,20).Value = 450
